# summer week 11 inputs
# Adds a new "Week 29" column (AD) with scores for the players that have
# reported a value for that week.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header for the newly added week
$ws.Range("AD1").Value = "Week 29"

# Scores for week 29, keyed by row (matching existing player rows)
$ws.Range("AD2").Value = 4     # Scott Foxley-Berry
$ws.Range("AD4").Value = 0     # Laura Thompson
$ws.Range("AD6").Value = 10    # Dan Aquino
$ws.Range("AD8").Value = 0     # Kim Quan
$ws.Range("AD9").Value = 2.5   # Leo Hayward

# Update the active selection to reflect where the user ended up working
$ws.Range("AE6").Select()
